$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 35: Quit Main Menu
$ws.Range("A35").Value = "Quit Main Menu"
$ws.Range("B35").Value = 430
$ws.Range("C35").Value = 668
$ws.Range("D35").Value = 158
$ws.Range("E35").Value = 70
$ws.Range("F35").Value = "Menu Box"

# New row 36: Back Credit Screen
$ws.Range("A36").Value = "Back Credit Screen"
$ws.Range("B36").Value = 428
$ws.Range("C36").Value = 668
$ws.Range("D36").Value = 198
$ws.Range("E36").Value = 68
$ws.Range("F36").Value = "Credit Screen"

# Copy style from row 33, whose per-column style pattern (s=2,5,5,5,5,2)
# matches the target rows exactly.
$ws.Range("A33:F33").Copy() | Out-Null
$ws.Range("A35:F35").PasteSpecial(-4122) | Out-Null
$ws.Range("A33:F33").Copy() | Out-Null
$ws.Range("A36:F36").PasteSpecial(-4122) | Out-Null

# Re-assert values since paste special (formats) shouldn't touch values, but ensure correctness
$ws.Range("A35").Value = "Quit Main Menu"
$ws.Range("B35").Value = 430
$ws.Range("C35").Value = 668
$ws.Range("D35").Value = 158
$ws.Range("E35").Value = 70
$ws.Range("F35").Value = "Menu Box"

$ws.Range("A36").Value = "Back Credit Screen"
$ws.Range("B36").Value = 428
$ws.Range("C36").Value = 668
$ws.Range("D36").Value = 198
$ws.Range("E36").Value = 68
$ws.Range("F36").Value = "Credit Screen"

$ws.Application.CutCopyMode = 0

# Update view
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C29").Select()
